# Auto-generated Excel COM-interop script to apply the diff changes
# Workbook sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 3500
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3500
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10500
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -10724

# ALC row 51
$ws.Range("H51").Value = 1499.875
$ws.Range("I51").Value = 1499
$ws.Range("K51").Value = 1499
$ws.Range("M51").Value = -1015

# ALC row 141
$ws.Range("H141").Value = 7998
$ws.Range("J141").Value = 7997
$ws.Range("L141").Value = 23991
$ws.Range("N141").Value = -34351

$ws = $wb.Worksheets.Item("ARM")
# ARM row 132
$ws.Range("H132").Value = 4398.143
$ws.Range("I132").Value = 2822.25
$ws.Range("K132").Value = 8466.75
$ws.Range("M132").Value = -5936.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 80
$ws.Range("H80").Value = 451.44446
$ws.Range("I80").Value = 123.333336
$ws.Range("K80").Value = 123.333336
$ws.Range("M80").Value = 874.666664

# BSM row 83
$ws.Range("H83").Value = 451.44446
$ws.Range("I83").Value = 123.333336
$ws.Range("K83").Value = 616.66668
$ws.Range("M83").Value = 4375.33332

# BSM row 141
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 159988
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 159988
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -170348

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 380
$ws.Range("I22").Value = 384.70587
$ws.Range("K22").Value = 384.70587
$ws.Range("M22").Value = -34.70587

# CRP row 86
$ws.Range("H86").Value = 20708.715
$ws.Range("I86").Value = 20708.715
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 20708.715
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -19585.715
$ws.Range("N86").ClearContents()

# CRP row 89
$ws.Range("H89").Value = 20708.715
$ws.Range("I89").Value = 20708.715
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 103543.575
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -97927.575
$ws.Range("N89").ClearContents()

# CRP row 99
$ws.Range("H99").Value = 2431
$ws.Range("I99").Value = 1574.6666
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 1574.6666
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -76.66660000000002
$ws.Range("N99").Value = -7996

# CRP row 126
$ws.Range("H126").Value = 2431
$ws.Range("I126").Value = 1574.6666
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 4723.9998
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -2253.9998
$ws.Range("N126").Value = -19940

# CRP row 132
$ws.Range("H132").Value = 4821.143
$ws.Range("J132").Value = 5249.5
$ws.Range("L132").Value = 15748.5
$ws.Range("N132").Value = -20808.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# CUL row 5
$ws.Range("H5").Value = 2663.3333
$ws.Range("I5").Value = 2663.3333
$ws.Range("K5").Value = 7989.999899999999
$ws.Range("M5").Value = -7877.999899999999

# CUL row 23
$ws.Range("H23").Value = 182.33333
$ws.Range("I23").Value = 23.5
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 70.5
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = 164.5
$ws.Range("N23").Value = -1970

# CUL row 29
$ws.Range("H29").Value = 1540.5
$ws.Range("J29").Value = 1720.5714
$ws.Range("L29").Value = 5161.7142
$ws.Range("N29").Value = -5715.7142

# CUL row 34
$ws.Range("H34").Value = 3101.1538
$ws.Range("J34").Value = 3833.4443
$ws.Range("L34").Value = 11500.3329
$ws.Range("N34").Value = -11668.3329

# CUL row 38
$ws.Range("H38").Value = 198
$ws.Range("I38").Value = 195
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 585
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = -238
$ws.Range("N38").Value = -1294

# CUL row 92
$ws.Range("H92").Value = 519.4
$ws.Range("I92").Value = 498
$ws.Range("J92").Value = 551.5
$ws.Range("K92").Value = 1494
$ws.Range("L92").Value = 1654.5
$ws.Range("M92").Value = -246
$ws.Range("N92").Value = -4150.5

# CUL row 135
$ws.Range("H135").Value = 2663.3333
$ws.Range("I135").Value = 2663.3333
$ws.Range("K135").Value = 23969.9997
$ws.Range("M135").Value = -21434.9997

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Range("H2").Value = 162.16667
$ws.Range("I2").Value = 194
$ws.Range("J2").Value = 66.666664
$ws.Range("K2").Value = 194
$ws.Range("L2").Value = 66.666664
$ws.Range("M2").Value = -81
$ws.Range("N2").Value = -292.666664

# GSM row 80
$ws.Range("H80").Value = 11428
$ws.Range("I80").Value = 4999
$ws.Range("J80").Value = 13999.6
$ws.Range("K80").Value = 4999
$ws.Range("L80").Value = 13999.6
$ws.Range("M80").Value = -4001
$ws.Range("N80").Value = -15995.6

# GSM row 83
$ws.Range("H83").Value = 11428
$ws.Range("I83").Value = 4999
$ws.Range("J83").Value = 13999.6
$ws.Range("K83").Value = 24995
$ws.Range("L83").Value = 69998
$ws.Range("M83").Value = -20003
$ws.Range("N83").Value = -79982

# GSM row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# GSM row 132
$ws.Range("H132").Value = 2802.7778
$ws.Range("I132").Value = 1782.8
$ws.Range("K132").Value = 5348.4
$ws.Range("M132").Value = -2818.4

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2500
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2724

# LTW row 16
$ws.Range("H16").Value = 1900
$ws.Range("I16").Value = 1900
$ws.Range("K16").Value = 1900
$ws.Range("M16").Value = -1730

# LTW row 68
$ws.Range("H68").Value = 3147.3333
$ws.Range("I68").Value = 1444
$ws.Range("J68").Value = 3999
$ws.Range("K68").Value = 1444
$ws.Range("L68").Value = 3999
$ws.Range("M68").Value = -695
$ws.Range("N68").Value = -5497

# LTW row 71
$ws.Range("H71").Value = 3147.3333
$ws.Range("I71").Value = 1444
$ws.Range("J71").Value = 3999
$ws.Range("K71").Value = 7220
$ws.Range("L71").Value = 19995
$ws.Range("M71").Value = -3476
$ws.Range("N71").Value = -27483

# LTW row 126
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 7500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -12440

# LTW row 132
$ws.Range("H132").Value = 6798.4165
$ws.Range("I132").Value = 1175.0555
$ws.Range("J132").Value = 23668.5
$ws.Range("K132").Value = 3525.1665
$ws.Range("L132").Value = 71005.5
$ws.Range("M132").Value = -995.1664999999998
$ws.Range("N132").Value = -76065.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 4000
$ws.Range("K96").Value = 4000
$ws.Range("M96").Value = -2627

# WVR row 126
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15440

# WVR row 132
$ws.Range("H132").Value = 2039.7
$ws.Range("I132").Value = 914.1429000000001
$ws.Range("K132").Value = 2742.4287
$ws.Range("M132").Value = -212.4287000000004
